$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# New rows of data to append (dates are Excel serial numbers: 44349-44353 = 2021-06-02..2021-06-06)
$data = @(
    @(44349, 11095, 264, 5075, 16434, 4380, 582, 25, 557, 113),
    @(44350, 11140, 107, 5139, 16486, 4435, 591, 25, 566, 113),
    @(44351, 11189, 211, 5184, 16584, 4471, 600, 25, 575, 113),
    @(44352, 11239, 187, 5232, 16658, 4531, 588, 25, 563, 113),
    @(44353, 11256, 167, 5236, 16659, 4581, 542, 25, 517, 113)
)

$startRow = 336
$lastExistingRow = 335

# Copy formatting from the last existing data row down onto the new rows,
# so they pick up the same cell styles (date format on col A, centered numbers elsewhere).
$srcRange = $ws.Range("A" + $lastExistingRow + ":J" + $lastExistingRow)
$srcRange.Copy()
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $destRange = $ws.Range("A" + $row + ":J" + $row)
    $destRange.PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    for ($col = 1; $col -le 10; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 1]
    }
}

$lastRow = $startRow + $data.Length - 1

$pane = $ws.Panes.Item(1)
$pane.TopLeftCell = $ws.Range("A325")
$ws.Range("B" + ($lastRow + 1)).Select()

$wb.Save()
